$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new row 36: schedule.reload ---
$ws.Rows.Item(36).Insert()
$ws.Range("A36").Value = "schedule.reload"
$ws.Range("B36").Value = "Làm mới toàn bộ dữ liệu"
$ws.Range("C36").Value = "Reload all data"
$ws.Rows.Item(36).RowHeight = 15
$ws.Range("A36").Font.Bold = $true

# --- Insert new row 94: scheduleRun.message.delete ---
# (after the first insertion, the old row 93 "scheduleRun.message.modifyStatus"
#  is now at row 93, so the new row goes at 94)
$ws.Rows.Item(94).Insert()
$ws.Range("A94").Value = "scheduleRun.message.delete"
$ws.Range("B94").Value = "Bạn có chắc chắn muốn xóa lịch trình chạy này?"
$ws.Range("C94").Value = "Do you want to delete this schedule run?"
$ws.Rows.Item(94).RowHeight = 15
$ws.Range("A94").Font.Bold = $true

# --- Update sheet view to match final selection state ---
$ws.Range("C95").Select()

Write-Output "done"
